$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'269.67"
$ws.Range("E2").Value = "'2.98%"
$ws.Range("G2").Value = "'21"
$ws.Range("D3").Value = "'26.71"
$ws.Range("E3").Value = "'-1.87%"
$ws.Range("G3").Value = "'21"
$ws.Range("E4").Value = "'0.07%"
$ws.Range("G4").Value = "'21"
$ws.Range("D5").Value = "'0.06102"
$ws.Range("E5").Value = "'-1.70%"
$ws.Range("G5").Value = "'21"
$ws.Range("D6").Value = "'6.738"
$ws.Range("E6").Value = "'0.29%"
$ws.Range("G6").Value = "'21"
$ws.Range("D7").Value = "'0.8586"
$ws.Range("E7").Value = "'1.00%"
$ws.Range("G7").Value = "'21"
$ws.Range("D8").Value = "'0.8913"
$ws.Range("E8").Value = "'-2.49%"
$ws.Range("G8").Value = "'21"
$ws.Range("D9").Value = "'0.1424"
$ws.Range("E9").Value = "'0.82%"
$ws.Range("G9").Value = "'21"
$ws.Range("D10").Value = "'0.05011"
$ws.Range("E10").Value = "'7.53%"
$ws.Range("G10").Value = "'21"
$ws.Range("D11").Value = "'0.07141"
$ws.Range("E11").Value = "'0.80%"
$ws.Range("G11").Value = "'21"
$ws.Range("D12").Value = "'0.03227"
$ws.Range("E12").Value = "'2.85%"
$ws.Range("G12").Value = "'21"
$ws.Range("D13").Value = "'0.09030"
$ws.Range("E13").Value = "'-0.28%"
$ws.Range("G13").Value = "'21"
$ws.Range("D14").Value = "'0.001545"
$ws.Range("E14").Value = "'0.14%"
$ws.Range("G14").Value = "'21"
$ws.Range("D15").Value = "'0.0006059"
$ws.Range("E15").Value = "'-1.99%"
$ws.Range("G15").Value = "'21"
$ws.Range("D16").Value = "'0.006083"
$ws.Range("E16").Value = "'1.01%"
$ws.Range("G16").Value = "'21"
$ws.Range("E17").Value = "'-0.09%"
$ws.Range("G17").Value = "'21"
$ws.Range("D18").Value = "'3.173"
$ws.Range("E18").Value = "'0.09%"
$ws.Range("G18").Value = "'21"
$ws.Range("E19").Value = "'2.91%"
$ws.Range("G19").Value = "'21"
$ws.Range("E20").Value = "'-0.63%"
$ws.Range("G20").Value = "'21"
$ws.Range("D21").Value = "'0.1299"
$ws.Range("E21").Value = "'-0.74%"
$ws.Range("G21").Value = "'21"
$ws.Range("D22").Value = "'3.840"
$ws.Range("E22").Value = "'-6.14%"
$ws.Range("G22").Value = "'21"
$ws.Range("D23").Value = "'0.04244"
$ws.Range("E23").Value = "'0.36%"
$ws.Range("G23").Value = "'21"
$ws.Range("D24").Value = "'0.001184"
$ws.Range("E24").Value = "'-2.04%"
$ws.Range("G24").Value = "'21"
$ws.Range("D25").Value = "'0.004147"
$ws.Range("E25").Value = "'9.07%"
$ws.Range("G25").Value = "'21"
$ws.Range("E26").Value = "'-0.03%"
$ws.Range("G26").Value = "'21"
$ws.Range("D27").Value = "'0.0001681"
$ws.Range("E27").Value = "'5.01%"
$ws.Range("G27").Value = "'21"
$ws.Range("G28").Value = "'21"
$ws.Range("G29").Value = "'21"
$ws.Range("G30").Value = "'21"
$ws.Range("G31").Value = "'21"
$ws.Range("G32").Value = "'21"
$ws.Range("G33").Value = "'21"
$ws.Range("G34").Value = "'21"
$ws.Range("G35").Value = "'21"
$ws.Range("G36").Value = "'21"
$ws.Range("G37").Value = "'21"
$ws.Range("G38").Value = "'21"
$ws.Range("G39").Value = "'21"
$ws.Range("D40").Value = "'0.03952"
$ws.Range("E40").Value = "'0.94%"
$ws.Range("G40").Value = "'21"
$ws.Range("D41").Value = "'0.1118"
$ws.Range("E41").Value = "'0.36%"
$ws.Range("G41").Value = "'21"
$ws.Range("D42").Value = "'0.004180"
$ws.Range("E42").Value = "'1.19%"
$ws.Range("G42").Value = "'21"
$ws.Range("D43").Value = "'0.002011"
$ws.Range("E43").Value = "'-7.93%"
$ws.Range("G43").Value = "'21"
$ws.Range("D44").Value = "'0.01273"
$ws.Range("E44").Value = "'-8.53%"
$ws.Range("G44").Value = "'21"
$ws.Range("D45").Value = "'0.00005134"
$ws.Range("E45").Value = "'-0.75%"
$ws.Range("G45").Value = "'21"
$ws.Range("E46").Value = "'-0.04%"
$ws.Range("G46").Value = "'21"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "'1.068"
$ws.Range("E47").Value = "'540.66%"
$ws.Range("G47").Value = "'21"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.02448"
$ws.Range("E48").Value = "'-31.82%"
$ws.Range("G48").Value = "'21"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("G49").Value = "'21"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("G50").Value = "'21"
$ws.Range("G51").Value = "'21"
Write-Output "OK"
